$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is a new log entry, duplicating most values from row 2 except the
# time (B), the comp. time (M) and — implicitly — that it's a distinct run.
#
# A3 ("2018.03.05") must stay as TEXT, matching the existing shared string
# used by A2. A direct Value assignment makes Excel auto-recognize the
# string as a date, so instead we pull the value in via a formula that
# references A2 (guaranteeing an exact text copy), then flatten the
# formula down to a plain cached value with a values-only paste.
$ws.Range("A3").Formula = "=A2"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$ws.Range("B3").Value = "14:24:28"
$ws.Range("C3").Value = "RS"
$ws.Range("D3").Value = 32
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = "N/A"
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = "effective"
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 5.98
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 35.4
$ws.Range("P3").Value = 0.5423719619281825
